$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# "Publish Java course materials on Anodiam platform" (row 18) is finished -> Status: Done
$ws.Range("F18").Value = "Done"

# "Publish AI course materials on Anodiam platform" (row 17) comment: note that the
# Python for AI course has now been uploaded to the Classplus website
$ws.Range("G17").Value = "PPTs only for now at version 0.0.1. Python for AI has been uploaded in Classplus website (Price given is Rs 1 as we will teach free)"

# Leave the cursor on G7, matching where editing left off
$ws.Range("G7").Select()

# Filter the tracker down to the items owned by Rahul
$ws.Range("A1:H25").AutoFilter(4, @("Rahul"), 7)
